$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("UserList")

# Update UserList (Sheet2) data: fix duplicate email and add a Flag column
$ws2.Range("B4").Value = "TestUSer2@gmail.com"
$ws2.Range("D2").Value = "NO"
$ws2.Range("D3").Value = "YES"
$ws2.Range("D4").Value = "NO"

# Update selections / active sheet to match the saved view state
$ws1.Range("F7").Select() | Out-Null
$ws2.Activate()
$ws2.Range("F15").Select() | Out-Null
